# Auto-generated: applies cell-level text updates per the OOXML diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds numeric-looking text (e.g. '126.80', '5.81').
# Writing such strings straight into .Value lets Excel auto-convert them to
# real floating point numbers (losing formatting / introducing FP noise).
# Force the whole column to Text format while we write, then restore the
# default style so no stray formatting is left behind.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '53.895.65'
$ws.Range("E2").Value = '  +0.78%  '
$ws.Range("D3").Value = '2.246.35'
$ws.Range("E3").Value = '  +2.44%  '
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("D5").Value = '491.94'
$ws.Range("E5").Value = '  +1.74%  '
$ws.Range("D6").Value = '126.80'
$ws.Range("E6").Value = '  +1.88%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  +0.27%  '
$ws.Range("D8").Value = '0.526'
$ws.Range("E8").Value = '  +1.31%  '
$ws.Range("D9").Value = '0.0950'
$ws.Range("E9").Value = '  +4.11%  '
$ws.Range("E10").Value = '  +2.53%  '
$ws.Range("E11").Value = '  +3.50%  '
$ws.Range("E12").Value = '  +1.06%  '
$ws.Range("D13").Value = '2.661.67'
$ws.Range("E13").Value = '  +3.11%  '
$ws.Range("D14").Value = '21.65'
$ws.Range("E14").Value = '  +3.15%  '
$ws.Range("D15").Value = '53.876.32'
$ws.Range("E15").Value = '  +0.89%  '
$ws.Range("D16").Value = '0.0000128'
$ws.Range("E16").Value = '  +1.04%  '
$ws.Range("D17").Value = '2.251.32'
$ws.Range("E17").Value = '  +2.61%  '
$ws.Range("D18").Value = '9.98'
$ws.Range("E18").Value = '  +4.67%  '
$ws.Range("D19").Value = '4.07'
$ws.Range("E19").Value = '  +3.49%  '
$ws.Range("D20").Value = '298.64'
$ws.Range("E20").Value = '  +1.41%  '
$ws.Range("D21").Value = '6.39'
$ws.Range("E21").Value = '  +5.72%  '
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  +0.07%  '
$ws.Range("E23").Value = '  -2.01%  '
$ws.Range("D24").Value = '61.79'
$ws.Range("E24").Value = '  -1.19%  '
$ws.Range("E25").Value = '  +2.30%  '
$ws.Range("D26").Value = '0.369'
$ws.Range("E26").Value = '  +1.26%  '
$ws.Range("D27").Value = '2.354.91'
$ws.Range("E27").Value = '  +2.66%  '
$ws.Range("D28").Value = '0.147'
$ws.Range("E28").Value = '  +2.10%  '
$ws.Range("D29").Value = '7.00'
$ws.Range("E29").Value = '  +0.77%  '
$ws.Range("D30").Value = '165.77'
$ws.Range("E30").Value = '  +0.22%  '
$ws.Range("E31").Value = '  +1.73%  '
$ws.Range("D32").Value = '0.0₃0678'
$ws.Range("E32").Value = '  +3.13%  '
$ws.Range("B33").Value = 'USDe'
$ws.Range("C33").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  +0.08%  '
$ws.Range("B34").Value = 'Aptos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D34").Value = '5.81'
$ws.Range("E34").Value = '  +2.49%  '
$ws.Range("D35").Value = '0.996'
$ws.Range("E35").Value = '  +0.21%  '
$ws.Range("E36").Value = '  +0.19%  '
$ws.Range("D37").Value = '17.58'
$ws.Range("E37").Value = '  +1.88%  '
$ws.Range("D38").Value = '0.896'
$ws.Range("E38").Value = '  +8.77%  '
$ws.Range("E39").Value = '  +2.79%  '
$ws.Range("D40").Value = '3.65'
$ws.Range("E40").Value = '  +3.59%  '
$ws.Range("D41").Value = '35.68'
$ws.Range("E41").Value = '  -0.03%  '
$ws.Range("D42").Value = '1.39'
$ws.Range("E42").Value = '  +2.52%  '
$ws.Range("D43").Value = '0.372'
$ws.Range("E43").Value = '  +1.79%  '
$ws.Range("D44").Value = '3.34'
$ws.Range("E44").Value = '  +3.31%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = '124.82'
$ws.Range("E45").Value = '  +1.43%  '
$ws.Range("B46").Value = 'RenderToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D46").Value = '4.71'
$ws.Range("E46").Value = '  -0.18%  '
$ws.Range("D47").Value = '0.0885'
$ws.Range("E47").Value = '  +1.45%  '
$ws.Range("D48").Value = '0.538'
$ws.Range("E48").Value = '  +1.52%  '
$ws.Range("D49").Value = '234.94'
$ws.Range("E49").Value = '  +2.90%  '
$ws.Range("D50").Value = '0.0481'
$ws.Range("E50").Value = '  +3.01%  '
$ws.Range("E51").Value = '  +1.69%  '

# Restore original (default) styling on the price column now that the
# text values are safely stored as strings.
$priceRange.Style = "Normal"
